# Updated with latest results.
# Fills in the measured/test values that were previously blank, flips the
# "Final decision" flags from "n" to "y", records the responsible tester,
# and stamps the test date. All the PASS/FAIL/ACCEPTED/REJECTED cells and
# the helper ok/NOK columns are formulas, so Excel recalculates them on
# its own once the raw inputs below are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: "Tests for shorts and connections" - K18:K33 observed values ---
$ws.Range("K18").Value = 0.7
$ws.Range("K19").Value = 0.7
$ws.Range("K20").Value = 0.7
$ws.Range("K21").Value = 0.7
$ws.Range("K22").Value = 0.7
$ws.Range("K23").Value = 0.7
$ws.Range("K24").Value = 0.7
$ws.Range("K25").Value = 0.7
$ws.Range("K26").Value = 0.7
$ws.Range("K27").Value = 0.7
$ws.Range("K28").Value = 0.7
$ws.Range("K29").Value = 0.7
$ws.Range("K30").Value = 0.7
$ws.Range("K31").Value = 0.7
$ws.Range("K32").Value = 0.7
$ws.Range("K33").Value = 0.7

# --- Section: dimension check resistances (K40:K47) ---
$ws.Range("K40").Value = 0.2
$ws.Range("K41").Value = 0.2
$ws.Range("K42").Value = 0.2
$ws.Range("K43").Value = 0.2
$ws.Range("K44").Value = 0.2
$ws.Range("K45").Value = 0.2
$ws.Range("K46").Value = 0.2
$ws.Range("K47").Value = 0.2

# --- Final decision flags: "n" -> "y" ---
$ws.Range("M57").Value = "y"
$ws.Range("M58").Value = "y"
$ws.Range("M59").Value = "y"

# Overall result of section 2 (LED test)
$ws.Range("B60").Value = "PASS"

# --- Section 3: HV test readings ---
$ws.Range("K66").Value = 98.8
$ws.Range("K67").Value = 9.92

# --- Section 4: voltage drop measurements ---
$ws.Range("C70").Value = 0.0246
$ws.Range("C71").Value = 0.034

# --- Tester initials and test date ---
$ws.Range("B80").Value = "Brian / Frank"
$ws.Range("F80").Value = "9/26/2014"
